# Updated Script for UI automation
# Renames the "ZipCode" column on the CheckOut sheet to "PostalCode" and
# formats the postal code value as an integer number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CheckOut")

# C1 header: "ZipCode" -> "PostalCode"
$ws.Range("C1").Value = "PostalCode"

# C2 holds the postal code value (214258); give it an integer number format
$ws.Range("C2").NumberFormat = "0"

# Re-fit column C now that its header/content changed
[void]$ws.Columns.Item(3).AutoFit()

# Move the active selection (as last left by the automation run)
[void]$ws.Range("G20").Select()
